$wb = $excel.ActiveWorkbook

# Rename the "VA" worksheet to "Product-VA"
$ws = $wb.Worksheets.Item("VA")
$ws.Name = "Product-VA"

# Update the scroll position / selection on the renamed sheet
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 231
$win.ScrollColumn = 1
$ws.Range("C257").Select()
